# Implementation of Cyl SDDR
#
# Adds a new "ITER Cylindrical benchmark for SDDR" entry (row 9) to the
# "Computational benchmarks" sheet, and makes that sheet the active tab
# (mirroring the author finishing their edit there instead of on
# "Experimental benchmarks").

$wb = $excel.ActiveWorkbook

$wsComp = $wb.Worksheets.Item("Computational benchmarks")
$wsExp  = $wb.Worksheets.Item("Experimental benchmarks")

# --- New benchmark row (row 9) on "Computational benchmarks" ---------------

# Plain text fields: ordinary .Value assignment is fine (not recognised as
# number/boolean/date, so it's stored as a normal shared string).
$wsComp.Range("A9").Value = "ITER Cylindrical benchmark for SDDR"
$wsComp.Range("B9").Value = "ITER_Cyl_SDDR.i"
$wsComp.Range("I9").Value = "D1S5"

# Run / OnlyInput / Post-Processing columns reuse the literal text "false"
# exactly like the row above (row 8). Typing "false" via .Value would be
# auto-coerced to a boolean, so instead copy the existing text cells from
# row 8 - this carries over both the literal string value and the plain
# (non quote-prefixed) cell style untouched.
$wsComp.Range("C8:E8").Copy()
$wsComp.Range("C9:E9").PasteSpecial()

# Numeric NPS cut-off, formatted like the other rows in this column
# (scientific notation).
$wsComp.Range("F9").Value = 500000000
$wsComp.Range("F9").NumberFormat = "0.00E+00"

# --- Selection / active-tab bookkeeping -------------------------------------

# Move the active selection to H9 on the Computational benchmarks sheet...
$wsComp.Range("H9").Select()

# ...and make "Computational benchmarks" the active sheet/tab instead of
# "Experimental benchmarks" (whose own selection, I5, is left untouched).
$wsComp.Activate()
